# Add new columns I (I0) and J (IF) to the sheet, matching the
# formatting used by the existing header cells (e.g. H1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, borders, centered/top aligned)
# from H1 onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-44
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 4
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 5
$ws.Range("I8").Value = 7
$ws.Range("J8").Value = 7
$ws.Range("I9").Value = 8
$ws.Range("J9").Value = 9
$ws.Range("I10").Value = 6
$ws.Range("J10").Value = 8
$ws.Range("I11").Value = 7
$ws.Range("J11").Value = 7
$ws.Range("I12").Value = 4
$ws.Range("J12").Value = 6
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 7
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 5
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 6
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 5
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 4
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 2
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = 6
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 7
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 5
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 5
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 6
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 3
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 5
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 6
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = 6
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = 4
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 5
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = 4
$ws.Range("I31").Value = 1
$ws.Range("J31").Value = 7
$ws.Range("I32").Value = 1
$ws.Range("J32").Value = 6
$ws.Range("I33").Value = 1
$ws.Range("J33").Value = 6
$ws.Range("I34").Value = 1
$ws.Range("J34").Value = 7
$ws.Range("I35").Value = 1
$ws.Range("J35").Value = 5
$ws.Range("I36").Value = 1
$ws.Range("J36").Value = 6
$ws.Range("I37").Value = 1
$ws.Range("J37").Value = 6
$ws.Range("I38").Value = 1
$ws.Range("J38").Value = 5
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = 5
$ws.Range("I40").Value = 1
$ws.Range("J40").Value = 3
$ws.Range("I41").Value = 1
$ws.Range("J41").Value = 5
$ws.Range("I42").Value = 1
$ws.Range("J42").Value = 4
$ws.Range("I43").Value = 1
$ws.Range("J43").Value = 3
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = 2
